$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fix the title wording, and the expectCode/expectMessage stay ---
$ws.Range("A2").Value = "图分析视图-queryLoadStatus-可正常获取图文件数据查询载入状态状态"

# --- Row 3: graphId missing case now expects 201 / 参数错误 / isRun=1 ---
$ws.Range("C3").Value = 201
$ws.Range("D3").Value = "参数错误"
$ws.Range("E3").Value = 1

# --- Row 4: projectId missing case, isRun now true ---
$ws.Range("E4").Value = 1

# --- Row 5: data source id missing case, isRun now true ---
$ws.Range("E5").Value = 1

# --- Row 6: no params case, isRun now true ---
$ws.Range("E6").Value = 1

# --- Row 7: new test case - cross-account cannot query load status ---
$ws.Range("A7").Value = "图分析视图-queryLoadStatus-跨账号不查询载入状态"
$ws.Range("B7").Value = '{"projectId":1334,"graphId":2497,"id":4113}'
$ws.Range("C7").Value = 401
$ws.Range("D7").Value = "无权访问"
$ws.Range("E7").Value = 1

# --- Row 8: new test case - same account, non current project graph cannot query ---
$ws.Range("A8").Value = "图分析视图-queryLoadStatus-同账号非当前项目标签不能查询载入状态"
$ws.Range("B8").Value = '{"projectId":1426,"graphId":2827,"id":4113}'
$ws.Range("C8").Value = 90009
$ws.Range("D8").Value = "无权操作"
$ws.Range("E8").Value = 1

# --- Row 9: new test case - graph not existing in project cannot query ---
$ws.Range("A9").Value = "图分析视图-queryLoadStatus-项目下不存在的标签不能查询载入状态"
$ws.Range("B9").Value = '{"projectId":1426,"graphId":2772,"id":4113}'
$ws.Range("C9").Value = 90009
$ws.Range("D9").Value = "无权操作"
$ws.Range("E9").Value = 1

# --- match the author's new active selection (single cell B7) ---
$ws.Range("B7").Select()
